$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column N: FmtBCDCol
$ws.Range("N1").Value = "FmtBCDCol"

$ws.Range("N2").Value = 2.2999999999999998

$ws.Range("N3").NumberFormat = "@"
$ws.Range("N3").Value = "200.24"
$ws.Range("N3").ClearFormats()

$ws.Range("N4").Value = -4000.99
$ws.Range("N5").Value = 200.91
$ws.Range("N6").Value = 1001
$ws.Range("N7").Value = 2002
$ws.Range("N8").Value = 1
$ws.Range("N9").Value = 3
$ws.Range("N10").Value = 5
$ws.Range("N11").Value = 9
$ws.Range("N12").Value = 1

$ws.Range("N11").Select()
